$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 10-17 (pitcher "Roblez") ---
$ws.Range("F10").Value = "FB"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Ball"

$ws.Range("F11").Value = "FB"
$ws.Range("G11").Value = "Take"
$ws.Range("H11").Value = "Ball"

$ws.Range("F12").Value = "FB"
$ws.Range("G12").Value = "Take"
$ws.Range("H12").Value = "Strike"
$ws.Range("M12").Value = $null

$ws.Range("F13").Value = "FB"
$ws.Range("G13").Value = "Swing"
$ws.Range("H13").Value = "Foul"

$ws.Range("F14").Value = "CH"
$ws.Range("G14").Value = "Take"
$ws.Range("H14").Value = "Ball"

$ws.Range("F15").Value = "FB"
$ws.Range("G15").Value = "Swing"
$ws.Range("H15").Value = "In Play"

$ws.Range("J17").Value = "CH,CB,FB"

# --- Block 2: rows 19-26 ---
$ws.Range("F19").Value = "FB"
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Ball"
$ws.Range("M19").Value = "83.5 MPH"

$ws.Range("F20").Value = "FB"
$ws.Range("G20").Value = "Take"
$ws.Range("H20").Value = "Ball"

$ws.Range("F21").Value = "FB"
$ws.Range("G21").Value = "Swing"
$ws.Range("H21").Value = "Foul"
$ws.Range("M21").Value = "32.75°"

$ws.Range("F22").Value = "FB"
$ws.Range("G22").Value = "Swing"
$ws.Range("H22").Value = "In Play"

$ws.Range("J26").Value = "CH,CB,FB"

# --- Block 3: rows 28-35 ---
$ws.Range("F28").Value = "CB"
$ws.Range("G28").Value = "Take"
$ws.Range("H28").Value = "Strike"
$ws.Range("M28").Value = "100.91 MPH"

$ws.Range("F29").Value = "CB"
$ws.Range("G29").Value = "Take"
$ws.Range("H29").Value = "Ball"

$ws.Range("F30").Value = "CH"
$ws.Range("G30").Value = "Take"
$ws.Range("H30").Value = "Ball"
$ws.Range("M30").Value = "8.06°"

$ws.Range("F31").Value = "FB"
$ws.Range("G31").Value = "Take"
$ws.Range("H31").Value = "Ball"

$ws.Range("F32").Value = "FB"
$ws.Range("G32").Value = "Swing"
$ws.Range("H32").Value = "In Play"

$ws.Range("J35").Value = "CH,CB,FB"

# --- Block 4: rows 37-44 ---
$ws.Range("F37").Value = "CH"
$ws.Range("G37").Value = "Take"
$ws.Range("H37").Value = "Ball"
$ws.Range("M37").Value = "82.59 MPH"

$ws.Range("F38").Value = "CH"
$ws.Range("G38").Value = "Take"
$ws.Range("H38").Value = "Strike"

$ws.Range("F39").Value = "CH"
$ws.Range("G39").Value = "Take"
$ws.Range("H39").Value = "Ball"
$ws.Range("M39").Value = "-15.88°"

$ws.Range("F40").Value = "CH"
$ws.Range("G40").Value = "Swing"
$ws.Range("H40").Value = "In Play"

$ws.Range("J44").Value = "CH,CB,FB,SL"
